# ValutazioneAnnotazioni.xlsx — "Aggiustati degli errori sul file name"
#
# The header labels get underscores instead of spaces, and several
# "punteggio equivalente" cells in column C that were mis-typed as
# decimal numbers (e.g. 20.23) get corrected to their intended
# semicolon-separated id-list text (e.g. "20; 23").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio1")

# --- Header row: rename labels (underscores instead of spaces). ---
# Set B1 ("punteggio") first so it lands as shared-string index 0,
# matching the target string table ordering.
$ws.Range("B1").Value = "punteggio"
$ws.Range("A1").Value = "id_annotatore"
$ws.Range("C1").Value = "id_annotatore_equivalente"

# --- Column C fixes: numbers that should have been id-lists. ---
$ws.Range("C18").Value = "20; 23"
$ws.Range("C20").Value = "22;26"
$ws.Range("C21").Value = "17;23"
$ws.Range("C23").Value = "19;26"
$ws.Range("C24").Value = "17;2"
$ws.Range("C27").Value = "19;22"

# --- Cosmetic: move the active selection like the author's session. ---
[void]$ws.Range("F17").Select()
